$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 159, pushing existing rows 159..202 down to 160..203
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row 159 with the new record's data
$ws.Cells.Item(159, 1).Value = 6
$ws.Cells.Item(159, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(159, 3).Value = "Metropolitana"
$ws.Cells.Item(159, 4).Value = 44627
$ws.Cells.Item(159, 5).Value = 13
$ws.Cells.Item(159, 6).Value = 100112022
$ws.Cells.Item(159, 7).Value = "Arveja Verde"
$ws.Cells.Item(159, 8).Value = "Sin especificar"
$ws.Cells.Item(159, 9).Value = "Primera"
$ws.Cells.Item(159, 10).Value = 280
$ws.Cells.Item(159, 11).Value = 23000
$ws.Cells.Item(159, 12).Value = 25000
$ws.Cells.Item(159, 13).Value = 23714
$ws.Cells.Item(159, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(159, 15).Value = "Carahue"
$ws.Cells.Item(159, 16).Value = 949
$ws.Cells.Item(159, 17).Value = 25
$ws.Cells.Item(159, 18).Value = "Hortaliza"
